$wb = $excel.ActiveWorkbook

# The "optimization_parameters" sheet had a stray leftover row (row 16:
# "Sheet", 3, 4) that needed to be removed while wrapping up the test
# file audit.
$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Activate() | Out-Null
$wsOpt.Rows.Item(16).Delete() | Out-Null

# The row that used to be 17 (the timepoints row) is now row 16; reselect
# it so the saved selection reflects the new layout.
$wsOpt.Range("A16:XFD16").Select() | Out-Null

# Wrap up by moving on to the next sheet (threshold_b) so it becomes the
# active / selected sheet when the workbook is saved.
$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate() | Out-Null
